# Fix wrong sheet name: "NewSheet" -> "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

# The author's last selection before saving was cell B10 (was C2).
$ws.Range("B10").Select() | Out-Null
